$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 11 ("References"): rework the "R ressources" section into
# "Unit testing" + extra reference links, and add a new "Debug mode" section
# with more links.
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$body = $s11.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: section heading "R ressources" -> "Unit testing"
$heading = $body.Paragraphs(1)
$heading.Text = "Unit testing"

# Paragraph 2: fix "RStudiom" typo -> "RStudio" (this naturally splits the
# run the same way PowerPoint would when only part of it is retyped).
$refPara = $body.Paragraphs(2)
$typo = $refPara.Characters(12, 10)
$typo.Text = "RStudio, "

# Insert the new paragraphs right after the reference paragraph, before the
# existing blank paragraphs that close out the placeholder. The leading `r
# breaks away from the reference paragraph before the new content starts.
$newText = "`r" + "https://fr.wikipedia.org/wiki/Test_unitaire" + "`r" + `
    "https://r-pkgs.org/tests.html" + "`r" + `
    "`r" + `
    "Debug mode" + "`r" + `
    "https://support.rstudio.com/hc/en-us/articles/205612627-Debugging-with-RStudio" + "`r" + `
    "`r"
$refPara.InsertAfter($newText)

# Paragraph 3: Wikipedia link on "unit test" (split the same way the
# original author typed it: "https" + "://" + rest).
$p3 = $body.Paragraphs(3)
$p3.Characters(1, 5).ActionSettings(1).Hyperlink.Address = "https://fr.wikipedia.org/wiki/Test_unitaire"
$p3.Characters(6, 3).ActionSettings(1).Hyperlink.Address = "https://fr.wikipedia.org/wiki/Test_unitaire"
$p3.Characters(9, 36).ActionSettings(1).Hyperlink.Address = "https://fr.wikipedia.org/wiki/Test_unitaire"

# Paragraph 4: r-pkgs.org testing chapter link.
$p4 = $body.Paragraphs(4)
$p4.Characters(1, 8).ActionSettings(1).Hyperlink.Address = "https://r-pkgs.org/tests.html"
$p4.Characters(9, 22).ActionSettings(1).Hyperlink.Address = "https://r-pkgs.org/tests.html"

# Paragraph 6: "Debug mode" sub-heading, bold + underlined like the other
# section heading.
$p6 = $body.Paragraphs(6)
$p6.Font.Bold = $true
$p6.Font.Underline = $true

# Paragraph 7: RStudio debugging support article link.
$p7 = $body.Paragraphs(7)
$p7.Characters(1, 8).ActionSettings(1).Hyperlink.Address = "https://support.rstudio.com/hc/en-us/articles/205612627-Debugging-with-RStudio"
$p7.Characters(9, 72).ActionSettings(1).Hyperlink.Address = "https://support.rstudio.com/hc/en-us/articles/205612627-Debugging-with-RStudio"

# ---------------------------------------------------------------------------
# Slide 11 footer date placeholders are driven by the layouts/master below;
# nothing shape-specific needed here.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Slide 5: merge a couple of runs that used to be split for no reason back
# into single runs (no visible text change, just cleanup).
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$testBody = $s5.Shapes.Item(2).TextFrame.TextRange
$forestPara = $testBody.Paragraphs(1)
$forestPara.Characters(1, 19).Text = "Forest model: test "

$popPara = $testBody.Paragraphs(2)
$popPara.Characters(65, 6).Text = " to 11"
$popPara.Characters(85, 22).Text = " are 6 " + [char]0x00AB + " mouse" + [char]0x00BB + " and 3 " + [char]0x00AB + " "

# ---------------------------------------------------------------------------
# Update the "updated on" date field cached text on every slide layout and
# on the slide master footer placeholder.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    foreach ($shp in $layout.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "4/27/2021") {
                $tr.Text = "6/7/2021"
            }
        }
    }
}

foreach ($shp in $p.SlideMaster.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "4/27/2021") {
            $tr.Text = "6/7/2021"
        }
    }
}
